$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.4135716542737669
$ws.Range("D2").Value = 0.03977780226842853
$ws.Range("E2").Value = 0.1747898442481883
$ws.Range("F2").Value = 1.060907389764225
$ws.Range("G2").Value = 0.002460226391477634
$ws.Range("K2").Value = 1.482885458775854
$ws.Range("L2").Value = 0.1526099929091203
$ws.Range("N2").Value = 1.20631993277059
$ws.Range("O2").Value = 3.776928602187127
$ws.Range("C3").Value = 0.4046509031364565
$ws.Range("D3").Value = 0.03917341654105755
$ws.Range("E3").Value = 0.1706213479563345
$ws.Range("F3").Value = 1.055964088884096
$ws.Range("G3").Value = 0.002463638155374274
$ws.Range("K3").Value = 1.341241579616224
$ws.Range("L3").Value = 0.1486158785716896
$ws.Range("N3").Value = 1.21327461594904
$ws.Range("O3").Value = 3.777273930213937
$ws.Range("C4").Value = 0.3994007138429367
$ws.Range("D4").Value = 0.03880172142994809
$ws.Range("E4").Value = 0.1681607113829884
$ws.Range("F4").Value = 1.05366492320082
$ws.Range("G4").Value = 0.00246584457850367
$ws.Range("K4").Value = 1.254430759288908
$ws.Range("L4").Value = 0.146252070831153
$ws.Range("N4").Value = 1.218034820363592
$ws.Range("O4").Value = 3.780042658595164
$ws.Range("C5").Value = 0.3973183445554014
$ws.Range("D5").Value = 0.03865011737331159
$ws.Range("E5").Value = 0.1671828182023773
$ws.Range("F5").Value = 1.052912730238205
$ws.Range("G5").Value = 0.002466771860239559
$ws.Range("K5").Value = 1.219096284854885
$ws.Range("L5").Value = 0.1453110479690665
$ws.Range("N5").Value = 1.220098023510722
$ws.Range("O5").Value = 3.781812600682031
$ws.Range("C6").Value = 0.3969760189206397
$ws.Range("D6").Value = 0.03862493592500371
$ws.Range("E6").Value = 0.1670219393460357
$ws.Range("F6").Value = 1.052798976620423
$ws.Range("G6").Value = 0.002466927537340428
$ws.Range("K6").Value = 1.213231585619724
$ws.Range("L6").Value = 0.1451561347060917
$ws.Range("N6").Value = 1.220448073776829
$ws.Range("O6").Value = 3.782145216955968
$ws.Range("C7").Value = 0.3993723989856051
$ws.Range("D7").Value = 0.03879967737128354
$ws.Range("E7").Value = 0.1681474226344477
$ws.Range("F7").Value = 1.053654031339612
$ws.Range("G7").Value = 0.002465856970065073
$ws.Range("K7").Value = 1.253954055439664
$ws.Range("L7").Value = 0.1462392898250187
$ws.Range("N7").Value = 1.218062145654869
$ws.Range("O7").Value = 3.780063932240495
$ws.Range("C8").Value = 0.4104486205837361
$ws.Range("D8").Value = 0.03956954364819154
$ws.Range("E8").Value = 0.1733320209080986
$ws.Range("F8").Value = 1.059049959863799
$ws.Range("G8").Value = 0.00246137966070813
$ws.Range("K8").Value = 1.434014562033326
$ws.Range("L8").Value = 0.1512144147799077
$ws.Range("N8").Value = 1.208616304877779
$ws.Range("O8").Value = 3.776516270554993
$ws.Range("C9").Value = 0.4339735444889357
$ws.Range("D9").Value = 0.04107388652258237
$ws.Range("E9").Value = 0.1842847787315662
$ws.Range("F9").Value = 1.075489443688397
$ws.Range("G9").Value = 0.002453481065173703
$ws.Range("K9").Value = 1.788321434045542
$ws.Range("L9").Value = 0.1616756405803983
$ws.Range("N9").Value = 1.193974362810067
$ws.Range("O9").Value = 3.789907073755273
$ws.Range("C10").Value = 0.4523624271112681
$ws.Range("D10").Value = 0.04217516513717356
$ws.Range("E10").Value = 0.1928140630590747
$ws.Range("F10").Value = 1.091166830522823
$ws.Range("G10").Value = 0.002448209692843182
$ws.Range("K10").Value = 2.049322490486702
$ws.Range("L10").Value = 0.169795191094579
$ws.Range("N10").Value = 1.185575226496212
$ws.Range("O10").Value = 3.812243133595445
$ws.Range("C11").Value = 0.4609692301995381
$ws.Range("D11").Value = 0.04267516828826246
$ws.Range("E11").Value = 0.1967997582766756
$ws.Range("F11").Value = 1.099086701078761
$ws.Range("G11").Value = 0.002445925876506268
$ws.Range("K11").Value = 2.168201103034107
$ws.Range("L11").Value = 0.1735840621492315
$ws.Range("N11").Value = 1.182264865233051
$ws.Range("O11").Value = 3.825138894416909
$ws.Range("C12").Value = 0.4642631970394859
$ws.Range("D12").Value = 0.04286435289379398
$ws.Range("E12").Value = 0.1983242743404787
$ws.Range("F12").Value = 1.102199533179842
$ws.Range("G12").Value = 0.002445077380163252
$ws.Range("K12").Value = 2.213237372207857
$ws.Range("L12").Value = 0.1750325576062295
$ws.Range("N12").Value = 1.181084606935272
$ws.Range("O12").Value = 3.830417017643043
$ws.Range("C13").Value = 0.4635522362652011
$ws.Range("D13").Value = 0.04282361577742222
$ws.Range("E13").Value = 0.1979952652545691
$ws.Range("F13").Value = 1.101524063460289
$ws.Range("G13").Value = 0.002445259393835752
$ws.Range("K13").Value = 2.203537165516479
$ws.Range("L13").Value = 0.1747199868733986
$ws.Range("N13").Value = 1.18133553795532
$ws.Range("O13").Value = 3.829262697646556
$ws.Range("C14").Value = 0.4612395301598724
$ws.Range("D14").Value = 0.04269073582936045
$ws.Range("E14").Value = 0.1969248759702751
$ws.Range("F14").Value = 1.09934051357186
$ws.Range("G14").Value = 0.002445855743311249
$ws.Range("K14").Value = 2.171905882816816
$ws.Range("L14").Value = 0.173702955321005
$ws.Range("N14").Value = 1.18216629609222
$ws.Range("O14").Value = 3.825565208270575
$ws.Range("C15").Value = 0.4598274577616337
$ws.Range("D15").Value = 0.0426093222814572
$ws.Range("E15").Value = 0.196271214382044
$ws.Range("F15").Value = 1.098017852013754
$ws.Range("G15").Value = 0.00244622314953523
$ws.Range("K15").Value = 2.152533298923117
$ws.Range("L15").Value = 0.1730817835400558
$ws.Range("N15").Value = 1.182684702961041
$ws.Range("O15").Value = 3.823351845484893
$ws.Range("C16").Value = 0.4518048152665415
$ws.Range("D16").Value = 0.04214246781114639
$ws.Range("E16").Value = 0.1925557170301104
$ws.Range("F16").Value = 1.090665143292185
$ws.Range("G16").Value = 0.002448361235706125
$ws.Range("K16").Value = 2.04155632007587
$ws.Range("L16").Value = 0.1695494991246562
$ws.Range("N16").Value = 1.185801829340903
$ws.Range("O16").Value = 3.811455527404036
$ws.Range("C17").Value = 0.4469450755017874
$ws.Range("D17").Value = 0.04185580702407776
$ws.Range("E17").Value = 0.1903034634068419
$ws.Range("F17").Value = 1.086356652641498
$ws.Range("G17").Value = 0.002449702063457575
$ws.Range("K17").Value = 1.973512173866993
$ws.Range("L17").Value = 0.1674069774077509
$ws.Range("N17").Value = 1.187844751082281
$ws.Range("O17").Value = 3.804859074032748
$ws.Range("C18").Value = 0.4441726221170939
$ws.Range("D18").Value = 0.04169083656431383
$ws.Range("E18").Value = 0.1890179724904044
$ws.Range("F18").Value = 1.083952682185455
$ws.Range("G18").Value = 0.002450484021951911
$ws.Range("K18").Value = 1.934389007854179
$ws.Range("L18").Value = 0.1661836184980814
$ws.Range("N18").Value = 1.189067836766398
$ws.Range("O18").Value = 3.80132230762635
$ws.Range("C19").Value = 0.4432378222992384
$ws.Range("D19").Value = 0.04163496529545085
$ws.Range("E19").Value = 0.1885844348600898
$ws.Range("F19").Value = 1.083151463456133
$ws.Range("G19").Value = 0.002450750628368007
$ws.Range("K19").Value = 1.921145038856139
$ws.Range("L19").Value = 0.165770948494611
$ws.Range("N19").Value = 1.189490208515082
$ws.Range("O19").Value = 3.800168972193973
$ws.Range("C20").Value = 0.4474600494512799
$ws.Range("D20").Value = 0.04188633206137382
$ws.Range("E20").Value = 0.1905421901330016
$ws.Range("F20").Value = 1.086807620606493
$ws.Range("G20").Value = 0.002449558218077441
$ws.Range("K20").Value = 1.980754150363794
$ws.Range("L20").Value = 0.1676341246168818
$ws.Range("N20").Value = 1.187622306195195
$ws.Range("O20").Value = 3.805534633318331
$ws.Range("C21").Value = 0.4619178845811689
$ws.Range("D21").Value = 0.04272977024177038
$ws.Range("E21").Value = 0.1972388620037009
$ws.Range("F21").Value = 1.099978784679791
$ws.Range("G21").Value = 0.002445680138306956
$ws.Range("K21").Value = 2.181196241476243
$ws.Range("L21").Value = 0.1740013091649786
$ws.Range("N21").Value = 1.181920293592256
$ws.Range("O21").Value = 3.826640524425045
$ws.Range("C22").Value = 0.4715695151158457
$ws.Range("D22").Value = 0.04328009098045982
$ws.Range("E22").Value = 0.201704247733872
$ws.Range("F22").Value = 1.109250111409338
$ws.Range("G22").Value = 0.00244324076640788
$ws.Range("K22").Value = 2.312309699959201
$ws.Range("L22").Value = 0.1782426911609747
$ws.Range("N22").Value = 1.178620928590007
$ws.Range("O22").Value = 3.842736213043622
$ws.Range("C23").Value = 0.4663997153979835
$ws.Range("D23").Value = 0.04298646327841738
$ws.Range("E23").Value = 0.1993128614101849
$ws.Range("F23").Value = 1.104241008863184
$ws.Range("G23").Value = 0.002444534022950618
$ws.Range("K23").Value = 2.242322219202038
$ws.Range("L23").Value = 0.175971648423328
$ws.Range("N23").Value = 1.180342800634833
$ws.Range("O23").Value = 3.833934534217207
$ws.Range("C24").Value = 0.447227162837919
$ws.Range("D24").Value = 0.04187253220501219
$ws.Range("E24").Value = 0.1904342326290305
$ws.Range("F24").Value = 1.086603510505682
$ws.Range("G24").Value = 0.002449623215997607
$ws.Range("K24").Value = 1.977480064120812
$ws.Range("L24").Value = 0.1675314051862671
$ws.Range("N24").Value = 1.187722722182244
$ws.Range("O24").Value = 3.805228416803658
$ws.Range("C25").Value = 0.4274157273099775
$ws.Range("D25").Value = 0.04066756847577579
$ws.Range("E25").Value = 0.181237309113925
$ws.Range("F25").Value = 1.070412086016191
$ws.Range("G25").Value = 0.002455524065573779
$ws.Range("K25").Value = 1.692347339377989
$ws.Range("L25").Value = 0.1587697240998125
$ws.Range("N25").Value = 1.197520716973052
$ws.Range("O25").Value = 3.784096522411318
